$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172146677970886
$ws.Range("B1").Value = 2.437282562255859
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.36595606803894
$ws.Range("E1").Value = 1.237607955932617
